$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F header - copy formatting (bold/border/center/top) from an
# existing header cell so it shares style index 1, then set its own value.
$ws.Range("F1").Value = "Trening"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Data for rows 2-13: Timestamp(serial), Seconds, Velocity, Acceleration_SMA, Velocity_Bin, Trening
$data = @(
    @(45686.4762449074,   1195.4, 10.2,  0.6324317455291748,  "10-15", "Duża Gra"),
    @(45686.47625648148,  1196.4, 10.87, 0.07279676584792982, "10-15", "Duża Gra"),
    @(45686.47626805556,  1197.4, 10.37, -0.2498368876320974, "10-15", "Duża Gra"),
    @(45686.4758212963,   1158.8, 7.42,  2.294765676770891,   "5-10",  "Duża Gra"),
    @(45686.47597986111,  1172.5, 6.17,  2.082028763634818,   "5-10",  "Duża Gra"),
    @(45686.47734791667,  1290.7, 6.62,  2.100573863301959,   "5-10",  "Duża Gra"),
    @(45686.48793819444,  2205.7, 14.74, 3.32965908731733,    "10-15", "Mała Gra"),
    @(45686.48802962963,  2213.6, 14.42, 3.529272590364729,   "10-15", "Mała Gra"),
    @(45686.49819050926,  3091.5, 10.93, 3.026102747235979,   "10-15", "Mała Gra"),
    @(45686.49735486111,  3019.3, 9.44,  2.834098151751926,   "5-10",  "Mała Gra"),
    @(45686.49802615741,  3077.3, 9.74,  2.9823728288923,     "5-10",  "Mała Gra"),
    @(45686.49818935185,  3091.4, 9.699999999999999, 2.861110533986773, "5-10", "Mała Gra")
)

# Set the number format twice (lowercase then uppercase) to reproduce the
# two numFmt entries (164 unused/lowercase, 165 used/uppercase) seen in the target.
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2:A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
}

Write-Output "done"
